$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values scraped verbatim from the source site
# (e.g. "28.202.64", "1.004") - not valid Excel numbers, so force Text
# formatting on the whole column before writing, otherwise Excel would
# coerce them to numeric and mangle formats like trailing zeros / double dots.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.223.49"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.797.54"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "331.64"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").Value = "0.4551"
$ws.Range("E7").Value = "  +17.36%  "
$ws.Range("D8").Value = "0.3735"
$ws.Range("E8").Value = "  +9.88%  "
$ws.Range("D9").Value = "44.68"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "1.144"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "0.07566"
$ws.Range("E11").Value = "  +4.91%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "22.41"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "6.313"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "7.515"
$ws.Range("E15").Value = "  +6.86%  "
$ws.Range("D16").Value = "1.790.96"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").Value = "0.00001092"
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("D18").Value = "0.06764"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("D19").Value = "80.94"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "17.53"
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("D22").Value = "6.353"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").Value = "28.209.42"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "11.82"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").Value = "2.430"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("D26").Value = "20.55"
$ws.Range("E26").Value = "  +3.13%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "152.05"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.357"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "1.994.14"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "133.38"
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("D31").Value = "1.238"
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("D32").Value = "4.037"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "0.09451"
$ws.Range("E33").Value = "  +9.57%  "
$ws.Range("D34").Value = "5.809"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").Value = "0.2368"
$ws.Range("E35").Value = "  +12.48%  "
$ws.Range("D36").Value = "12.13"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "0.06329"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("D38").Value = "0.02336"
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").Value = "5.202"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("D40").Value = "0.6587"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").Value = "8.371"
$ws.Range("E41").Value = "  +6.64%  "
$ws.Range("D42").Value = "1.483"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("D43").Value = "1.207"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "14.22"
$ws.Range("E44").Value = "  +4.00%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "0.6112"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").Value = "3.805"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Value = "129.96"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("D49").Value = "2.033"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").Value = "0.07124"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").Value = "1.161"
$ws.Range("E51").Value = "  +0.61%  "
